# Update scripts with new tpm values (Vtn-Tnfrsf11b LR-pair data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update recalculated TPM-derived values in existing rows (2-7) ---
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.337313
$ws.Cells.Item(2, 8).Value = 10.011939
$ws.Cells.Item(2, 9).Value = 0.1958858017947999
$ws.Cells.Item(2, 10).Value = 0.1958858017947999
$ws.Cells.Item(2, 13).Value = 0.08241233333333334
$ws.Cells.Item(2, 17).Value = 0.2750357513936667
$ws.Cells.Item(2, 18).Value = 2.475321762543
$ws.Cells.Item(2, 19).Value = 0.009259532845638938
$ws.Cells.Item(2, 20).Value = 0.009259532845638938
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.337313
$ws.Cells.Item(3, 8).Value = 10.011939
$ws.Cells.Item(3, 9).Value = 0.1958858017947999
$ws.Cells.Item(3, 10).Value = 0.1958858017947999
$ws.Cells.Item(3, 17).Value = 5.543356988512
$ws.Cells.Item(3, 18).Value = 49.890212896608
$ws.Cells.Item(3, 19).Value = 0.186626268949161
$ws.Cells.Item(3, 20).Value = 0.186626268949161
$ws.Cells.Item(4, 9).Value = 0.6036083824529627
$ws.Cells.Item(4, 10).Value = 0.6036083824529627
$ws.Cells.Item(4, 13).Value = 0.08241233333333334
$ws.Cells.Item(4, 17).Value = 0.8475034101214445
$ws.Cells.Item(4, 18).Value = 7.627530691093001
$ws.Cells.Item(4, 19).Value = 0.02853260211825403
$ws.Cells.Item(4, 20).Value = 0.02853260211825403
$ws.Cells.Item(5, 9).Value = 0.6036083824529627
$ws.Cells.Item(5, 10).Value = 0.6036083824529627
$ws.Cells.Item(5, 19).Value = 0.5750757803347086
$ws.Cells.Item(5, 20).Value = 0.5750757803347086
$ws.Cells.Item(6, 7).Value = 3.355061
$ws.Cells.Item(6, 8).Value = 10.065183
$ws.Cells.Item(6, 9).Value = 0.196927532435664
$ws.Cells.Item(6, 10).Value = 0.196927532435664
$ws.Cells.Item(6, 13).Value = 0.08241233333333334
$ws.Cells.Item(6, 17).Value = 0.2764984054856667
$ws.Cells.Item(6, 18).Value = 2.488485649371
$ws.Cells.Item(6, 19).Value = 0.009308775511503482
$ws.Cells.Item(6, 20).Value = 0.009308775511503482
$ws.Cells.Item(7, 7).Value = 3.355061
$ws.Cells.Item(7, 8).Value = 10.065183
$ws.Cells.Item(7, 9).Value = 0.196927532435664
$ws.Cells.Item(7, 10).Value = 0.196927532435664
$ws.Cells.Item(7, 17).Value = 5.572836842464
$ws.Cells.Item(7, 18).Value = 50.155531582176
$ws.Cells.Item(7, 19).Value = 0.1876187569241606
$ws.Cells.Item(7, 20).Value = 0.1876187569241605

# --- Append two new rows (8-9) for the "Resolving-Mac" sending cluster ---
$ws.Cells.Item(8, 1).Value = "Resolving-Mac"
$ws.Cells.Item(8, 2).Value = "Vtn"
$ws.Cells.Item(8, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.06096333333333333
$ws.Cells.Item(8, 8).Value = 0.18289
$ws.Cells.Item(8, 9).Value = 0.003578283316573439
$ws.Cells.Item(8, 10).Value = 0.003578283316573439
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.08241233333333334
$ws.Cells.Item(8, 14).Value = 0.247237
$ws.Cells.Item(8, 15).Value = 0.04727005612861496
$ws.Cells.Item(8, 16).Value = 0.04727005612861496
$ws.Cells.Item(8, 17).Value = 0.005024130547777778
$ws.Cells.Item(8, 18).Value = 0.04521717493
$ws.Cells.Item(8, 19).Value = 0.000169145653218513
$ws.Cells.Item(8, 20).Value = 0.000169145653218513
$ws.Cells.Item(9, 1).Value = "Resolving-Mac"
$ws.Cells.Item(9, 2).Value = "Vtn"
$ws.Cells.Item(9, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.06096333333333333
$ws.Cells.Item(9, 8).Value = 0.18289
$ws.Cells.Item(9, 9).Value = 0.003578283316573439
$ws.Cells.Item(9, 10).Value = 0.003578283316573439
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.661024
$ws.Cells.Item(9, 14).Value = 4.983072
$ws.Cells.Item(9, 15).Value = 0.9527299438713851
$ws.Cells.Item(9, 16).Value = 0.952729943871385
$ws.Cells.Item(9, 17).Value = 0.1012615597866667
$ws.Cells.Item(9, 18).Value = 0.9113540380799999
$ws.Cells.Item(9, 19).Value = 0.003409137663354927
$ws.Cells.Item(9, 20).Value = 0.003409137663354926

Write-Output ("Dimension now: " + $ws.UsedRange.Address())
